$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.998.55'
$ws.Range("E2").Value = '  +1.26%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.957.81'
$ws.Range("E3").Value = '  -0.13%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.69'
$ws.Range("E5").Value = '  -1.30%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4863'
$ws.Range("E7").Value = '  +1.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2951'
$ws.Range("E8").Value = '  +0.99%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06933'
$ws.Range("E9").Value = '  +2.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.45'
$ws.Range("E10").Value = '  +1.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '107.86'
$ws.Range("E11").Value = '  -1.13%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.003.74'
$ws.Range("E12").Value = '  +2.15%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07804'
$ws.Range("E13").Value = '  +1.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.502'
$ws.Range("E14").Value = '  +1.30%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7017'
$ws.Range("E15").Value = '  +1.34%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '281.92'
$ws.Range("E16").Value = '  -3.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.037.97'
$ws.Range("E17").Value = '  +1.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.31'
$ws.Range("E18").Value = '  +1.29%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007782'
$ws.Range("E19").Value = '  +1.51%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.210.04'
$ws.Range("E20").Value = '  +0.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.28%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.522'
$ws.Range("E22").Value = '  -2.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9989'
$ws.Range("E23").Value = '  -0.19%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.523'
$ws.Range("E24").Value = '  -1.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.882'
$ws.Range("E25").Value = '  -0.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.72'
$ws.Range("E26").Value = '  -0.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.96'
$ws.Range("E27").Value = '  +0.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.195'
$ws.Range("E28").Value = '  +1.04%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1050'
$ws.Range("E29").Value = '  -1.36%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.394'
$ws.Range("E30").Value = '  -3.14%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.663'
$ws.Range("E31").Value = '  -2.46%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.573'
$ws.Range("E32").Value = '  -1.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.467'
$ws.Range("E33").Value = '  +0.70%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04937'
$ws.Range("E34").Value = '  -2.80%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7560'
$ws.Range("E35").Value = '  -1.63%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.171'
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.732'
$ws.Range("E37").Value = '  +0.44%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02011'
$ws.Range("E38").Value = '  -0.82%  '

# Row 39
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.571'
$ws.Range("E40").Value = '  +1.59%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.57'
$ws.Range("E41").Value = '  +12.48%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.124'
$ws.Range("E42").Value = '  -0.71%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9046'
$ws.Range("E43").Value = '  +2.53%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.61'
$ws.Range("E44").Value = '  -0.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4455'
$ws.Range("E45").Value = '  +0.11%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.184'
$ws.Range("E46").Value = '  +9.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007.06'
$ws.Range("E48").Value = '  +8.96%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.361'
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1256'
$ws.Range("E50").Value = '  -1.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.03'
$ws.Range("E51").Value = '  +0.34%  '
